$d = $word.ActiveDocument
$null = $d.Content.Find.Execute("11+54=", $true, $false, $false, $false, $false, $true, 1, $false, "48-24=", 2)
$null = $d.Content.Find.Execute("39-39=", $true, $false, $false, $false, $false, $true, 1, $false, "8+38=", 2)
$null = $d.Content.Find.Execute("53+14=", $true, $false, $false, $false, $false, $true, 1, $false, "82-13=", 2)
$null = $d.Content.Find.Execute("96-74=", $true, $false, $false, $false, $false, $true, 1, $false, "75-62=", 2)
$null = $d.Content.Find.Execute("34+30=", $true, $false, $false, $false, $false, $true, 1, $false, "26-7=", 2)
$null = $d.Content.Find.Execute("22+13=", $true, $false, $false, $false, $false, $true, 1, $false, "27+48=", 2)
$null = $d.Content.Find.Execute("24+15=", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=", 2)
$null = $d.Content.Find.Execute("18-6=", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=", 2)
$null = $d.Content.Find.Execute("48-8=", $true, $false, $false, $false, $false, $true, 1, $false, "84-6=", 2)
$null = $d.Content.Find.Execute("28+6=", $true, $false, $false, $false, $false, $true, 1, $false, "12+82=", 2)
$null = $d.Content.Find.Execute("53-31=", $true, $false, $false, $false, $false, $true, 1, $false, "85-39=", 2)
$null = $d.Content.Find.Execute("30-22=", $true, $false, $false, $false, $false, $true, 1, $false, "29+62=", 2)
$null = $d.Content.Find.Execute("79+1=", $true, $false, $false, $false, $false, $true, 1, $false, "83-72=", 2)
$null = $d.Content.Find.Execute("30+24=", $true, $false, $false, $false, $false, $true, 1, $false, "72-42=", 2)
$null = $d.Content.Find.Execute("87-41=", $true, $false, $false, $false, $false, $true, 1, $false, "59-38=", 2)
$null = $d.Content.Find.Execute("40+28=", $true, $false, $false, $false, $false, $true, 1, $false, "70-54=", 2)
$null = $d.Content.Find.Execute("88-37=", $true, $false, $false, $false, $false, $true, 1, $false, "37-1=", 2)
$null = $d.Content.Find.Execute("60-54=", $true, $false, $false, $false, $false, $true, 1, $false, "60-15=", 2)
$null = $d.Content.Find.Execute("45+1=", $true, $false, $false, $false, $false, $true, 1, $false, "57-44=", 2)
$null = $d.Content.Find.Execute("99-20=", $true, $false, $false, $false, $false, $true, 1, $false, "55-6=", 2)
$null = $d.Content.Find.Execute("52+28=", $true, $false, $false, $false, $false, $true, 1, $false, "65-60=", 2)
$null = $d.Content.Find.Execute("21+52=", $true, $false, $false, $false, $false, $true, 1, $false, "68+6=", 2)
$null = $d.Content.Find.Execute("21+74=", $true, $false, $false, $false, $false, $true, 1, $false, "36+62=", 2)
$null = $d.Content.Find.Execute("28+50=", $true, $false, $false, $false, $false, $true, 1, $false, "89-20=", 2)
$null = $d.Content.Find.Execute("28-11=", $true, $false, $false, $false, $false, $true, 1, $false, "92-86=", 2)
$null = $d.Content.Find.Execute("58-47=", $true, $false, $false, $false, $false, $true, 1, $false, "20+65=", 2)
$null = $d.Content.Find.Execute("54-52=", $true, $false, $false, $false, $false, $true, 1, $false, "78-40=", 2)
$null = $d.Content.Find.Execute("16+37=", $true, $false, $false, $false, $false, $true, 1, $false, "23-12=", 2)
$null = $d.Content.Find.Execute("61-22=", $true, $false, $false, $false, $false, $true, 1, $false, "37-36=", 2)
$null = $d.Content.Find.Execute("29+33=", $true, $false, $false, $false, $false, $true, 1, $false, "9+23=", 2)
$null = $d.Content.Find.Execute("55-15=", $true, $false, $false, $false, $false, $true, 1, $false, "7+3=", 2)
$null = $d.Content.Find.Execute("13+66=", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=", 2)
$null = $d.Content.Find.Execute("14+72=", $true, $false, $false, $false, $false, $true, 1, $false, "55-14=", 2)
$null = $d.Content.Find.Execute("54-27=", $true, $false, $false, $false, $false, $true, 1, $false, "52-34=", 2)
$null = $d.Content.Find.Execute("35+21=", $true, $false, $false, $false, $false, $true, 1, $false, "22+1=", 2)
$null = $d.Content.Find.Execute("75-25=", $true, $false, $false, $false, $false, $true, 1, $false, "76-50=", 2)
$null = $d.Content.Find.Execute("1+40=", $true, $false, $false, $false, $false, $true, 1, $false, "20-19=", 2)
$null = $d.Content.Find.Execute("61-21=", $true, $false, $false, $false, $false, $true, 1, $false, "80+0=", 2)
$null = $d.Content.Find.Execute("69-8=", $true, $false, $false, $false, $false, $true, 1, $false, "83-48=", 2)
$null = $d.Content.Find.Execute("0+22=", $true, $false, $false, $false, $false, $true, 1, $false, "28+9=", 2)
$null = $d.Content.Find.Execute("59+24=", $true, $false, $false, $false, $false, $true, 1, $false, "10-4=", 2)
$null = $d.Content.Find.Execute("27-7=", $true, $false, $false, $false, $false, $true, 1, $false, "64-37=", 2)
$null = $d.Content.Find.Execute("90-42=", $true, $false, $false, $false, $false, $true, 1, $false, "60-18=", 2)
$null = $d.Content.Find.Execute("55+11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+79=", 2)
$null = $d.Content.Find.Execute("74-48=", $true, $false, $false, $false, $false, $true, 1, $false, "0+31=", 2)
$null = $d.Content.Find.Execute("43-41=", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=", 2)
$null = $d.Content.Find.Execute("9+10=", $true, $false, $false, $false, $false, $true, 1, $false, "73-64=", 2)
$null = $d.Content.Find.Execute("42+33=", $true, $false, $false, $false, $false, $true, 1, $false, "62-22=", 2)
$null = $d.Content.Find.Execute("81+9=", $true, $false, $false, $false, $false, $true, 1, $false, "28+8=", 2)
$null = $d.Content.Find.Execute("62-49=", $true, $false, $false, $false, $false, $true, 1, $false, "11+68=", 2)
$null = $d.Content.Find.Execute("29+65=", $true, $false, $false, $false, $false, $true, 1, $false, "70-50=", 2)
$null = $d.Content.Find.Execute("45-3=", $true, $false, $false, $false, $false, $true, 1, $false, "37+3=", 2)
$null = $d.Content.Find.Execute("10+82=", $true, $false, $false, $false, $false, $true, 1, $false, "42+19=", 2)
$null = $d.Content.Find.Execute("93-63=", $true, $false, $false, $false, $false, $true, 1, $false, "25-11=", 2)
$null = $d.Content.Find.Execute("48+31=", $true, $false, $false, $false, $false, $true, 1, $false, "96-31=", 2)
$null = $d.Content.Find.Execute("71-26=", $true, $false, $false, $false, $false, $true, 1, $false, "47+9=", 2)
$null = $d.Content.Find.Execute("24+27=", $true, $false, $false, $false, $false, $true, 1, $false, "26+6=", 2)
$null = $d.Content.Find.Execute("0+1=", $true, $false, $false, $false, $false, $true, 1, $false, "24-20=", 2)
$null = $d.Content.Find.Execute("49+23=", $true, $false, $false, $false, $false, $true, 1, $false, "41-27=", 2)
$null = $d.Content.Find.Execute("51+24=", $true, $false, $false, $false, $false, $true, 1, $false, "77-4=", 2)
$null = $d.Content.Find.Execute("35-15=", $true, $false, $false, $false, $false, $true, 1, $false, "17+76=", 2)
$null = $d.Content.Find.Execute("1+52=", $true, $false, $false, $false, $false, $true, 1, $false, "54+45=", 2)
$null = $d.Content.Find.Execute("69-39=", $true, $false, $false, $false, $false, $true, 1, $false, "19+53=", 2)
$null = $d.Content.Find.Execute("20+55=", $true, $false, $false, $false, $false, $true, 1, $false, "8+75=", 2)
$null = $d.Content.Find.Execute("20+43=", $true, $false, $false, $false, $false, $true, 1, $false, "92-39=", 2)
$null = $d.Content.Find.Execute("47+23=", $true, $false, $false, $false, $false, $true, 1, $false, "76-40=", 2)
$null = $d.Content.Find.Execute("99-66=", $true, $false, $false, $false, $false, $true, 1, $false, "53-47=", 2)
$null = $d.Content.Find.Execute("14+85=", $true, $false, $false, $false, $false, $true, 1, $false, "43+34=", 2)
$null = $d.Content.Find.Execute("47-28=", $true, $false, $false, $false, $false, $true, 1, $false, "57+19=", 2)
$null = $d.Content.Find.Execute("42+25=", $true, $false, $false, $false, $false, $true, 1, $false, "59-23=", 2)
$null = $d.Content.Find.Execute("20+18=", $true, $false, $false, $false, $false, $true, 1, $false, "70-21=", 2)
$null = $d.Content.Find.Execute("94-90=", $true, $false, $false, $false, $false, $true, 1, $false, "93-14=", 2)
$null = $d.Content.Find.Execute("71+15=", $true, $false, $false, $false, $false, $true, 1, $false, "30+6=", 2)
$null = $d.Content.Find.Execute("48-25=", $true, $false, $false, $false, $false, $true, 1, $false, "10+35=", 2)
$null = $d.Content.Find.Execute("46-30=", $true, $false, $false, $false, $false, $true, 1, $false, "26+73=", 2)
$null = $d.Content.Find.Execute("4+74=", $true, $false, $false, $false, $false, $true, 1, $false, "18+15=", 2)
$null = $d.Content.Find.Execute("31+4=", $true, $false, $false, $false, $false, $true, 1, $false, "66-22=", 2)
$null = $d.Content.Find.Execute("75+15=", $true, $false, $false, $false, $false, $true, 1, $false, "16+41=", 2)
$null = $d.Content.Find.Execute("60-36=", $true, $false, $false, $false, $false, $true, 1, $false, "8+62=", 2)
$null = $d.Content.Find.Execute("58+27=", $true, $false, $false, $false, $false, $true, 1, $false, "6+3=", 2)
$null = $d.Content.Find.Execute("48+22=", $true, $false, $false, $false, $false, $true, 1, $false, "27+9=", 2)
$null = $d.Content.Find.Execute("6+24=", $true, $false, $false, $false, $false, $true, 1, $false, "27-21=", 2)
$null = $d.Content.Find.Execute("23+41=", $true, $false, $false, $false, $false, $true, 1, $false, "55-42=", 2)
$null = $d.Content.Find.Execute("82-3=", $true, $false, $false, $false, $false, $true, 1, $false, "14-7=", 2)
$null = $d.Content.Find.Execute("97-19=", $true, $false, $false, $false, $false, $true, 1, $false, "94-32=", 2)
$null = $d.Content.Find.Execute("98-28=", $true, $false, $false, $false, $false, $true, 1, $false, "91-5=", 2)
$null = $d.Content.Find.Execute("37+58=", $true, $false, $false, $false, $false, $true, 1, $false, "6+89=", 2)
$null = $d.Content.Find.Execute("37-34=", $true, $false, $false, $false, $false, $true, 1, $false, "80-54=", 2)
$null = $d.Content.Find.Execute("47-15=", $true, $false, $false, $false, $false, $true, 1, $false, "20-0=", 2)
$null = $d.Content.Find.Execute("97-23=", $true, $false, $false, $false, $false, $true, 1, $false, "93-24=", 2)
$null = $d.Content.Find.Execute("1+32=", $true, $false, $false, $false, $false, $true, 1, $false, "64-8=", 2)
$null = $d.Content.Find.Execute("47+10=", $true, $false, $false, $false, $false, $true, 1, $false, "0+14=", 2)
$null = $d.Content.Find.Execute("17+35=", $true, $false, $false, $false, $false, $true, 1, $false, "36+33=", 2)
$null = $d.Content.Find.Execute("51-25=", $true, $false, $false, $false, $false, $true, 1, $false, "34+25=", 2)
$null = $d.Content.Find.Execute("21+42=", $true, $false, $false, $false, $false, $true, 1, $false, "4+5=", 2)
$null = $d.Content.Find.Execute("47+41=", $true, $false, $false, $false, $false, $true, 1, $false, "25+19=", 2)
$null = $d.Content.Find.Execute("48-26=", $true, $false, $false, $false, $false, $true, 1, $false, "0+51=", 2)
$null = $d.Content.Find.Execute("39-16=", $true, $false, $false, $false, $false, $true, 1, $false, "2+91=", 2)
$null = $d.Content.Find.Execute("11+21=", $true, $false, $false, $false, $false, $true, 1, $false, "54+22=", 2)
$null = $d.Content.Find.Execute("74-22=", $true, $false, $false, $false, $false, $true, 1, $false, "59-4=", 2)
